$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 25157.334
$ws.Range("I98").Value = 26190.2
$ws.Range("J98").Value = 4500
$ws.Range("K98").Value = 26190.2
$ws.Range("L98").Value = 4500
$ws.Range("M98").Value = -24692.2
$ws.Range("N98").Value = -7496
# Row 113
$ws.Range("H113").Value = 303186.7
$ws.Range("I113").Value = 430127.2
$ws.Range("J113").Value = 3145.4546
$ws.Range("K113").Value = 430127.2
$ws.Range("L113").Value = 3145.4546
$ws.Range("M113").Value = -426873.2
$ws.Range("N113").Value = -9653.454600000001
# Row 122
$ws.Range("H122").Value = 25157.334
$ws.Range("I122").Value = 26190.2
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 78570.60000000001
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -76120.60000000001
$ws.Range("N122").Value = -18400
# Row 137
$ws.Range("H137").Value = 864.3077
$ws.Range("I137").Value = 853.25
$ws.Range("J137").Value = 997
$ws.Range("K137").Value = 2559.75
$ws.Range("L137").Value = 2991
$ws.Range("M137").Value = -9.75
$ws.Range("N137").Value = -8091

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1412.25
$ws.Range("I2").Value = 1241
$ws.Range("K2").Value = 1241
$ws.Range("M2").Value = -1128
# Row 61
$ws.Range("H61").Value = 1421.3334
$ws.Range("I61").Value = 1410.1333
$ws.Range("J61").Value = 1533.3334
$ws.Range("K61").Value = 1410.1333
$ws.Range("L61").Value = 1533.3334
$ws.Range("M61").Value = -1198.1333
$ws.Range("N61").Value = -1957.3334
# Row 74
$ws.Range("H74").Value = 834.8182
$ws.Range("I74").Value = 676.8333
$ws.Range("J74").Value = 1024.4
$ws.Range("K74").Value = 676.8333
$ws.Range("L74").Value = 1024.4
$ws.Range("M74").Value = 197.1667
$ws.Range("N74").Value = -2772.4
# Row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# Row 77
$ws.Range("H77").Value = 834.8182
$ws.Range("I77").Value = 676.8333
$ws.Range("J77").Value = 1024.4
$ws.Range("K77").Value = 3384.1665
$ws.Range("L77").Value = 5122
$ws.Range("M77").Value = 983.8334999999997
$ws.Range("N77").Value = -13858
# Row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
# Row 94
$ws.Range("H94").Value = 32999.5
$ws.Range("J94").Value = 32999.5
$ws.Range("L94").Value = 32999.5
$ws.Range("N94").Value = -34801.5
# Row 116
$ws.Range("H116").Value = 1412.25
$ws.Range("I116").Value = 1241
$ws.Range("K116").Value = 1241
$ws.Range("M116").Value = 1053
# Row 136
$ws.Range("H136").Value = 1421.3334
$ws.Range("I136").Value = 1410.1333
$ws.Range("J136").Value = 1533.3334
$ws.Range("K136").Value = 4230.3999
$ws.Range("L136").Value = 4600.0002
$ws.Range("M136").Value = -1680.3999
$ws.Range("N136").Value = -9700.0002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1412.25
$ws.Range("I3").Value = 1241
$ws.Range("K3").Value = 1241
$ws.Range("M3").Value = -1127
# Row 94
$ws.Range("H94").Value = 329
$ws.Range("I94").Value = 330.9
$ws.Range("J94").Value = 310
$ws.Range("K94").Value = 330.9
$ws.Range("L94").Value = 310
$ws.Range("M94").Value = 120.1
$ws.Range("N94").Value = -1212
# Row 134
$ws.Range("H134").Value = 1614.9231
$ws.Range("I134").Value = 1028.75
$ws.Range("K134").Value = 3086.25
$ws.Range("M134").Value = -551.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2331.5454
$ws.Range("I31").Value = 2069.75
$ws.Range("J31").Value = 3029.6667
$ws.Range("K31").Value = 2069.75
$ws.Range("L31").Value = 3029.6667
$ws.Range("M31").Value = -1774.75
$ws.Range("N31").Value = -3619.6667
# Row 34
$ws.Range("H34").Value = 2331.5454
$ws.Range("I34").Value = 2069.75
$ws.Range("J34").Value = 3029.6667
$ws.Range("K34").Value = 2069.75
$ws.Range("L34").Value = 3029.6667
$ws.Range("M34").Value = -1867.75
$ws.Range("N34").Value = -3433.6667
# Row 58
$ws.Range("H58").Value = 929.05884
$ws.Range("I58").Value = 741.0741
$ws.Range("J58").Value = 1654.1428
$ws.Range("K58").Value = 741.0741
$ws.Range("L58").Value = 1654.1428
$ws.Range("M58").Value = -538.0741
$ws.Range("N58").Value = -2060.1428
# Row 136
$ws.Range("H136").Value = 929.05884
$ws.Range("I136").Value = 741.0741
$ws.Range("J136").Value = 1654.1428
$ws.Range("K136").Value = 2223.2223
$ws.Range("L136").Value = 4962.428400000001
$ws.Range("M136").Value = 326.7776999999996
$ws.Range("N136").Value = -10062.4284

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 22499.285
$ws.Range("J68").Value = 12997.667
$ws.Range("L68").Value = 38993.001
$ws.Range("N68").Value = -40615.001
# Row 69
$ws.Range("H69").Value = 1192.25
$ws.Range("I69").Value = 967.8333
$ws.Range("J69").Value = 1416.6666
$ws.Range("K69").Value = 2903.4999
$ws.Range("L69").Value = 4249.9998
$ws.Range("M69").Value = -2092.4999
$ws.Range("N69").Value = -5871.9998
# Row 71
$ws.Range("H71").Value = 22499.285
$ws.Range("J71").Value = 12997.667
$ws.Range("L71").Value = 116979.003
$ws.Range("N71").Value = -125091.003
# Row 72
$ws.Range("H72").Value = 1192.25
$ws.Range("I72").Value = 967.8333
$ws.Range("J72").Value = 1416.6666
$ws.Range("K72").Value = 8710.4997
$ws.Range("L72").Value = 12749.9994
$ws.Range("M72").Value = -4654.4997
$ws.Range("N72").Value = -20861.9994
# Row 86
$ws.Range("I86").Value = 500
$ws.Range("K86").Value = 1500
$ws.Range("M86").Value = -314
# Row 89
$ws.Range("I89").Value = 500
$ws.Range("K89").Value = 4500
$ws.Range("M89").Value = 1428

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2856.4614
$ws.Range("I102").Value = 2157.3635
$ws.Range("J102").Value = 3761.1765
$ws.Range("K102").Value = 2157.3635
$ws.Range("L102").Value = 3761.1765
$ws.Range("M102").Value = -535.3634999999999
$ws.Range("N102").Value = -7005.1765
# Row 132
$ws.Range("H132").Value = 5323.433
$ws.Range("I132").Value = 5511.2964
$ws.Range("J132").Value = 3632.6667
$ws.Range("K132").Value = 16533.8892
$ws.Range("L132").Value = 10898.0001
$ws.Range("M132").Value = -14003.8892
$ws.Range("N132").Value = -15958.0001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 391.94116
$ws.Range("I22").Value = 397.36365
$ws.Range("J22").Value = 382
$ws.Range("K22").Value = 397.36365
$ws.Range("L22").Value = 382
$ws.Range("M22").Value = -102.36365
$ws.Range("N22").Value = -972
# Row 27
$ws.Range("H27").Value = 391.94116
$ws.Range("I27").Value = 397.36365
$ws.Range("J27").Value = 382
$ws.Range("K27").Value = 397.36365
$ws.Range("L27").Value = 382
$ws.Range("M27").Value = -290.36365
$ws.Range("N27").Value = -596
# Row 68
$ws.Range("H68").Value = 7532.5
$ws.Range("I68").Value = 8755.385
$ws.Range("J68").Value = 2233.3333
$ws.Range("K68").Value = 8755.385
$ws.Range("L68").Value = 2233.3333
$ws.Range("M68").Value = -8006.385
$ws.Range("N68").Value = -3731.3333
# Row 71
$ws.Range("H71").Value = 7532.5
$ws.Range("I71").Value = 8755.385
$ws.Range("J71").Value = 2233.3333
$ws.Range("K71").Value = 43776.925
$ws.Range("L71").Value = 11166.6665
$ws.Range("M71").Value = -40032.925
$ws.Range("N71").Value = -18654.6665
# Row 76
$ws.Range("H76").Value = 5761
$ws.Range("I76").Value = 5761
$ws.Range("K76").Value = 5761
$ws.Range("M76").Value = -5423
# Row 79
$ws.Range("H79").Value = 5761
$ws.Range("I79").Value = 5761
$ws.Range("K79").Value = 5761
$ws.Range("M79").Value = -4591
# Row 122
$ws.Range("H122").Value = 3357.5789
$ws.Range("I122").Value = 3234.9412
$ws.Range("J122").Value = 4400
$ws.Range("K122").Value = 9704.8236
$ws.Range("L122").Value = 13200
$ws.Range("M122").Value = -7254.8236
$ws.Range("N122").Value = -18100

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 475.45
$ws.Range("I113").Value = 363.5
$ws.Range("J113").Value = 550.0833
$ws.Range("K113").Value = 1090.5
$ws.Range("L113").Value = 1650.2499
$ws.Range("M113").Value = 1079.5
$ws.Range("N113").Value = -5990.2499
# Row 126
$ws.Range("H126").Value = 1261
$ws.Range("I126").Value = 1034.4
$ws.Range("K126").Value = 3103.2
$ws.Range("M126").Value = -633.2000000000003
